$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"
$newOverviewDate = "2016-11-23 12:00:04"
$newZhHandoffDate = "2016-11-23 11:59:50"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/40227450588301bc56ff78991e553066bf2b3601/e2e/2c37e513-929f-40fd-ad3a-7fac559c4422.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/70c8d01e1cd22392d5e86e6bffeff1da702c098e/e2e/2c37e513-929f-40fd-ad3a-7fac559c4422.md."

# --- Overview sheet ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $newZhHandoffDate
$wsZhCn.Range("P2").Value = $errorDetail

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $newOverviewDate
$wsDeDe.Range("P2").Value = $errorDetail

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667
